$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 12).Value = 'stimuli/img_mgnmm.png'
$ws.Cells.Item(2, 13).Value = 79.14705882352941
$ws.Cells.Item(2, 14).Value = 60.38235294117647
$ws.Cells.Item(2, 15).Value = 69.76470588235294
$ws.Cells.Item(2, 16).Value = 34
$ws.Cells.Item(2, 17).Value = 8
$ws.Cells.Item(2, 18).Value = 8
$ws.Cells.Item(2, 19).Value = 8
$ws.Cells.Item(2, 20).Value = 8
$ws.Cells.Item(2, 21).Value = 8
$ws.Cells.Item(2, 22).Value = 8

# Row 3
$ws.Cells.Item(3, 9).ClearContents()
$ws.Cells.Item(3, 10).Value = 'new'
$ws.Cells.Item(3, 11).Value = 'f'
$ws.Cells.Item(3, 12).Value = 'stimuli/img_wyl6z.png'
$ws.Cells.Item(3, 13).Value = 59.8235294117647
$ws.Cells.Item(3, 14).Value = 36.23529411764706
$ws.Cells.Item(3, 15).Value = 48.02941176470588
$ws.Cells.Item(3, 16).Value = 34
$ws.Cells.Item(3, 17).Value = 3
$ws.Cells.Item(3, 18).Value = 3
$ws.Cells.Item(3, 19).Value = 3
$ws.Cells.Item(3, 20).Value = 3
$ws.Cells.Item(3, 21).Value = 3
$ws.Cells.Item(3, 22).Value = 3

# Row 4
$ws.Cells.Item(4, 12).Value = 'stimuli/img_ce9vx.png'
$ws.Cells.Item(4, 13).Value = 75.90909090909091
$ws.Cells.Item(4, 14).Value = 57.12121212121212
$ws.Cells.Item(4, 15).Value = 66.51515151515152
$ws.Cells.Item(4, 17).Value = 7
$ws.Cells.Item(4, 18).Value = 7
$ws.Cells.Item(4, 19).Value = 7
$ws.Cells.Item(4, 20).Value = 7
$ws.Cells.Item(4, 21).Value = 7
$ws.Cells.Item(4, 22).Value = 7

# Row 5
$ws.Cells.Item(5, 9).ClearContents()
$ws.Cells.Item(5, 10).Value = 'new'
$ws.Cells.Item(5, 11).Value = 'f'
$ws.Cells.Item(5, 12).Value = 'stimuli/img_nyv2b.png'
$ws.Cells.Item(5, 13).Value = 11.91176470588235
$ws.Cells.Item(5, 14).Value = 6.852941176470588
$ws.Cells.Item(5, 15).Value = 9.382352941176471
$ws.Cells.Item(5, 16).Value = 34
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = 1
$ws.Cells.Item(5, 19).Value = 1
$ws.Cells.Item(5, 20).Value = 1
$ws.Cells.Item(5, 21).Value = 1
$ws.Cells.Item(5, 22).Value = 1

# Row 6
$ws.Cells.Item(6, 9).ClearContents()
$ws.Cells.Item(6, 10).Value = 'new'
$ws.Cells.Item(6, 11).Value = 'f'
$ws.Cells.Item(6, 12).Value = 'stimuli/img_30vz5.png'
$ws.Cells.Item(6, 13).Value = 86.21212121212122
$ws.Cells.Item(6, 14).Value = 68.27272727272727
$ws.Cells.Item(6, 15).Value = 77.24242424242425
$ws.Cells.Item(6, 16).Value = 33
$ws.Cells.Item(6, 17).Value = 10
$ws.Cells.Item(6, 18).Value = 10
$ws.Cells.Item(6, 19).Value = 10
$ws.Cells.Item(6, 20).Value = 10
$ws.Cells.Item(6, 21).Value = 10
$ws.Cells.Item(6, 22).Value = 10

# Row 7
$ws.Cells.Item(7, 9).Value = 'target'
$ws.Cells.Item(7, 10).Value = 'old'
$ws.Cells.Item(7, 11).Value = 'j'
$ws.Cells.Item(7, 12).Value = 'stimuli/img_esb4r.png'
$ws.Cells.Item(7, 13).Value = 60.73529411764706
$ws.Cells.Item(7, 14).Value = 38.58823529411764
$ws.Cells.Item(7, 15).Value = 49.66176470588235
$ws.Cells.Item(7, 16).Value = 34
$ws.Cells.Item(7, 17).Value = 3
$ws.Cells.Item(7, 18).Value = 3
$ws.Cells.Item(7, 19).Value = 3
$ws.Cells.Item(7, 20).Value = 3
$ws.Cells.Item(7, 21).Value = 3
$ws.Cells.Item(7, 22).Value = 3

# Row 8
$ws.Cells.Item(8, 9).Value = 'target'
$ws.Cells.Item(8, 10).Value = 'old'
$ws.Cells.Item(8, 11).Value = 'j'
$ws.Cells.Item(8, 12).Value = 'stimuli/img_1ao2d.png'
$ws.Cells.Item(8, 13).Value = 38.77777777777778
$ws.Cells.Item(8, 14).Value = 18.75
$ws.Cells.Item(8, 15).Value = 28.76388888888889
$ws.Cells.Item(8, 16).Value = 36
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = 1
$ws.Cells.Item(8, 19).Value = 1
$ws.Cells.Item(8, 20).Value = 1
$ws.Cells.Item(8, 21).Value = 1
$ws.Cells.Item(8, 22).Value = 1

# Row 9
$ws.Cells.Item(9, 9).Value = 'target'
$ws.Cells.Item(9, 10).Value = 'old'
$ws.Cells.Item(9, 11).Value = 'j'
$ws.Cells.Item(9, 12).Value = 'stimuli/img_wppku.png'
$ws.Cells.Item(9, 13).Value = 75.02941176470588
$ws.Cells.Item(9, 14).Value = 53.05882352941177
$ws.Cells.Item(9, 15).Value = 64.04411764705883
$ws.Cells.Item(9, 16).Value = 34
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = 6
$ws.Cells.Item(9, 19).Value = 6
$ws.Cells.Item(9, 20).Value = 6
$ws.Cells.Item(9, 21).Value = 6
$ws.Cells.Item(9, 22).Value = 6

# Row 10
$ws.Cells.Item(10, 12).Value = 'stimuli/img_es7o2.png'
$ws.Cells.Item(10, 13).Value = 52.48571428571429
$ws.Cells.Item(10, 14).Value = 27.54285714285714
$ws.Cells.Item(10, 15).Value = 40.01428571428572
$ws.Cells.Item(10, 16).Value = 35
$ws.Cells.Item(10, 17).Value = 2
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 19).Value = 2
$ws.Cells.Item(10, 20).Value = 2
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 2

# Row 11
$ws.Cells.Item(11, 9).Value = 'target'
$ws.Cells.Item(11, 10).Value = 'old'
$ws.Cells.Item(11, 11).Value = 'j'
$ws.Cells.Item(11, 12).Value = 'stimuli/img_7ucnr.png'
$ws.Cells.Item(11, 13).Value = 70.39393939393939
$ws.Cells.Item(11, 14).Value = 47.90909090909091
$ws.Cells.Item(11, 15).Value = 59.15151515151515
$ws.Cells.Item(11, 17).Value = 5
$ws.Cells.Item(11, 18).Value = 5
$ws.Cells.Item(11, 19).Value = 5
$ws.Cells.Item(11, 20).Value = 5
$ws.Cells.Item(11, 21).Value = 5
$ws.Cells.Item(11, 22).Value = 5

# Row 12
$ws.Cells.Item(12, 9).ClearContents()
$ws.Cells.Item(12, 10).Value = 'new'
$ws.Cells.Item(12, 11).Value = 'f'
$ws.Cells.Item(12, 12).Value = 'stimuli/img_inqod.png'
$ws.Cells.Item(12, 13).Value = 70.84848484848484
$ws.Cells.Item(12, 14).Value = 50.63636363636363
$ws.Cells.Item(12, 15).Value = 60.74242424242424
$ws.Cells.Item(12, 16).Value = 33
$ws.Cells.Item(12, 17).Value = 5
$ws.Cells.Item(12, 18).Value = 5
$ws.Cells.Item(12, 19).Value = 5
$ws.Cells.Item(12, 20).Value = 5
$ws.Cells.Item(12, 21).Value = 5
$ws.Cells.Item(12, 22).Value = 5

# Row 13
$ws.Cells.Item(13, 9).ClearContents()
$ws.Cells.Item(13, 10).Value = 'new'
$ws.Cells.Item(13, 11).Value = 'f'
$ws.Cells.Item(13, 12).Value = 'stimuli/img_ye5sl.png'
$ws.Cells.Item(13, 13).Value = 53.2258064516129
$ws.Cells.Item(13, 14).Value = 34.45161290322581
$ws.Cells.Item(13, 15).Value = 43.83870967741936
$ws.Cells.Item(13, 16).Value = 31
$ws.Cells.Item(13, 17).Value = 2
$ws.Cells.Item(13, 18).Value = 2
$ws.Cells.Item(13, 19).Value = 2
$ws.Cells.Item(13, 20).Value = 2
$ws.Cells.Item(13, 21).Value = 2
$ws.Cells.Item(13, 22).Value = 2

# Row 14
$ws.Cells.Item(14, 12).Value = 'stimuli/img_6nbgt.png'
$ws.Cells.Item(14, 13).Value = 78.45161290322581
$ws.Cells.Item(14, 14).Value = 57.83870967741935
$ws.Cells.Item(14, 15).Value = 68.14516129032258
$ws.Cells.Item(14, 16).Value = 31
$ws.Cells.Item(14, 17).Value = 7
$ws.Cells.Item(14, 18).Value = 7
$ws.Cells.Item(14, 19).Value = 7
$ws.Cells.Item(14, 20).Value = 7
$ws.Cells.Item(14, 21).Value = 7
$ws.Cells.Item(14, 22).Value = 7

# Row 15
$ws.Cells.Item(15, 12).Value = 'stimuli/img_ikk62.png'
$ws.Cells.Item(15, 13).Value = 37.48780487804878
$ws.Cells.Item(15, 14).Value = 21.07317073170732
$ws.Cells.Item(15, 15).Value = 29.28048780487805
$ws.Cells.Item(15, 16).Value = 41
$ws.Cells.Item(15, 17).Value = 1
$ws.Cells.Item(15, 18).Value = 1
$ws.Cells.Item(15, 19).Value = 1
$ws.Cells.Item(15, 20).Value = 1
$ws.Cells.Item(15, 21).Value = 1
$ws.Cells.Item(15, 22).Value = 1

# Row 16
$ws.Cells.Item(16, 12).Value = 'stimuli/img_p3hpc.png'
$ws.Cells.Item(16, 13).Value = 72.83333333333333
$ws.Cells.Item(16, 14).Value = 52.22222222222222
$ws.Cells.Item(16, 15).Value = 62.52777777777777
$ws.Cells.Item(16, 16).Value = 36
$ws.Cells.Item(16, 17).Value = 6
$ws.Cells.Item(16, 18).Value = 6
$ws.Cells.Item(16, 19).Value = 6
$ws.Cells.Item(16, 20).Value = 6
$ws.Cells.Item(16, 21).Value = 6
$ws.Cells.Item(16, 22).Value = 6

# Row 17
$ws.Cells.Item(17, 9).Value = 'target'
$ws.Cells.Item(17, 10).Value = 'old'
$ws.Cells.Item(17, 11).Value = 'j'
$ws.Cells.Item(17, 12).Value = 'stimuli/img_mawe6.png'
$ws.Cells.Item(17, 13).Value = 83.48387096774194
$ws.Cells.Item(17, 14).Value = 65.54838709677419
$ws.Cells.Item(17, 15).Value = 74.51612903225806
$ws.Cells.Item(17, 16).Value = 31
$ws.Cells.Item(17, 17).Value = 9
$ws.Cells.Item(17, 18).Value = 9
$ws.Cells.Item(17, 19).Value = 9
$ws.Cells.Item(17, 20).Value = 9
$ws.Cells.Item(17, 21).Value = 9
$ws.Cells.Item(17, 22).Value = 9

# Row 18
$ws.Cells.Item(18, 9).ClearContents()
$ws.Cells.Item(18, 10).Value = 'new'
$ws.Cells.Item(18, 11).Value = 'f'
$ws.Cells.Item(18, 12).Value = 'stimuli/img_cnyac.png'
$ws.Cells.Item(18, 13).Value = 69.14705882352941
$ws.Cells.Item(18, 14).Value = 47.8235294117647
$ws.Cells.Item(18, 15).Value = 58.48529411764706
$ws.Cells.Item(18, 16).Value = 34
$ws.Cells.Item(18, 17).Value = 5
$ws.Cells.Item(18, 18).Value = 5
$ws.Cells.Item(18, 19).Value = 5
$ws.Cells.Item(18, 20).Value = 5
$ws.Cells.Item(18, 21).Value = 5
$ws.Cells.Item(18, 22).Value = 5

# Row 19
$ws.Cells.Item(19, 9).ClearContents()
$ws.Cells.Item(19, 10).Value = 'new'
$ws.Cells.Item(19, 11).Value = 'f'
$ws.Cells.Item(19, 12).Value = 'stimuli/img_d8xbu.png'
$ws.Cells.Item(19, 13).Value = 91.36363636363636
$ws.Cells.Item(19, 14).Value = 73.18181818181819
$ws.Cells.Item(19, 15).Value = 82.27272727272728
$ws.Cells.Item(19, 16).Value = 33
$ws.Cells.Item(19, 17).Value = 10
$ws.Cells.Item(19, 18).Value = 10
$ws.Cells.Item(19, 19).Value = 10
$ws.Cells.Item(19, 20).Value = 10
$ws.Cells.Item(19, 21).Value = 10
$ws.Cells.Item(19, 22).Value = 10

# Row 20
$ws.Cells.Item(20, 12).Value = 'stimuli/img_qmgwq.png'
$ws.Cells.Item(20, 13).Value = 84.58333333333333
$ws.Cells.Item(20, 14).Value = 64.44444444444444
$ws.Cells.Item(20, 15).Value = 74.51388888888889
$ws.Cells.Item(20, 16).Value = 36
$ws.Cells.Item(20, 17).Value = 9
$ws.Cells.Item(20, 18).Value = 9
$ws.Cells.Item(20, 19).Value = 9
$ws.Cells.Item(20, 20).Value = 9
$ws.Cells.Item(20, 21).Value = 9
$ws.Cells.Item(20, 22).Value = 9

# Row 21
$ws.Cells.Item(21, 12).Value = 'stimuli/img_kwxq1.png'
$ws.Cells.Item(21, 13).Value = 68.53125
$ws.Cells.Item(21, 14).Value = 44.09375
$ws.Cells.Item(21, 15).Value = 56.3125
$ws.Cells.Item(21, 16).Value = 32
$ws.Cells.Item(21, 17).Value = 4
$ws.Cells.Item(21, 18).Value = 4
$ws.Cells.Item(21, 19).Value = 4
$ws.Cells.Item(21, 20).Value = 4
$ws.Cells.Item(21, 21).Value = 4
$ws.Cells.Item(21, 22).Value = 4

# Row 22
$ws.Cells.Item(22, 12).Value = 'stimuli/img_zi8qc.png'
$ws.Cells.Item(22, 13).Value = 77.14285714285714
$ws.Cells.Item(22, 14).Value = 57.02857142857143
$ws.Cells.Item(22, 15).Value = 67.08571428571429
$ws.Cells.Item(22, 16).Value = 35
$ws.Cells.Item(22, 17).Value = 7
$ws.Cells.Item(22, 18).Value = 7
$ws.Cells.Item(22, 19).Value = 7
$ws.Cells.Item(22, 20).Value = 7
$ws.Cells.Item(22, 21).Value = 7
$ws.Cells.Item(22, 22).Value = 7

# Row 23
$ws.Cells.Item(23, 12).Value = 'stimuli/img_eatdk.png'
$ws.Cells.Item(23, 13).Value = 81.40625
$ws.Cells.Item(23, 14).Value = 61.375
$ws.Cells.Item(23, 15).Value = 71.390625
$ws.Cells.Item(23, 16).Value = 32
$ws.Cells.Item(23, 17).Value = 8
$ws.Cells.Item(23, 18).Value = 8
$ws.Cells.Item(23, 19).Value = 8
$ws.Cells.Item(23, 20).Value = 8
$ws.Cells.Item(23, 21).Value = 8
$ws.Cells.Item(23, 22).Value = 8

# Row 25
$ws.Cells.Item(25, 8).Value = 'kitchens'
$ws.Cells.Item(25, 10).Value = 'new'
$ws.Cells.Item(25, 12).Value = 'stimuli/img_aplao.png'
$ws.Cells.Item(25, 13).Value = 64.09090909090909
$ws.Cells.Item(25, 14).Value = 40.75757575757576
$ws.Cells.Item(25, 15).Value = 52.42424242424242
$ws.Cells.Item(25, 16).Value = 33
$ws.Cells.Item(25, 17).Value = 3
$ws.Cells.Item(25, 18).Value = 3
$ws.Cells.Item(25, 19).Value = 3
$ws.Cells.Item(25, 20).Value = 3
$ws.Cells.Item(25, 21).Value = 3
$ws.Cells.Item(25, 22).Value = 3

# Row 26
$ws.Cells.Item(26, 9).Value = 'target'
$ws.Cells.Item(26, 10).Value = 'old'
$ws.Cells.Item(26, 11).Value = 'j'
$ws.Cells.Item(26, 12).Value = 'stimuli/img_89rmb.png'
$ws.Cells.Item(26, 13).Value = 55.18518518518518
$ws.Cells.Item(26, 14).Value = 29.25925925925926
$ws.Cells.Item(26, 15).Value = 42.22222222222222
$ws.Cells.Item(26, 16).Value = 27
$ws.Cells.Item(26, 17).Value = 2
$ws.Cells.Item(26, 18).Value = 2
$ws.Cells.Item(26, 19).Value = 2
$ws.Cells.Item(26, 20).Value = 2
$ws.Cells.Item(26, 21).Value = 2
$ws.Cells.Item(26, 22).Value = 2

# Row 27
$ws.Cells.Item(27, 12).Value = 'stimuli/img_xdhz2.png'
$ws.Cells.Item(27, 13).Value = 63.3
$ws.Cells.Item(27, 14).Value = 37.25
$ws.Cells.Item(27, 15).Value = 50.275
$ws.Cells.Item(27, 16).Value = 40
$ws.Cells.Item(27, 17).Value = 3
$ws.Cells.Item(27, 18).Value = 3
$ws.Cells.Item(27, 19).Value = 3
$ws.Cells.Item(27, 20).Value = 3
$ws.Cells.Item(27, 21).Value = 3
$ws.Cells.Item(27, 22).Value = 3

# Row 28
$ws.Cells.Item(28, 12).Value = 'stimuli/img_7w5tw.png'
$ws.Cells.Item(28, 13).Value = 53.2258064516129
$ws.Cells.Item(28, 14).Value = 28.90322580645161
$ws.Cells.Item(28, 15).Value = 41.06451612903226
$ws.Cells.Item(28, 16).Value = 31
$ws.Cells.Item(28, 17).Value = 2
$ws.Cells.Item(28, 18).Value = 2
$ws.Cells.Item(28, 19).Value = 2
$ws.Cells.Item(28, 20).Value = 2
$ws.Cells.Item(28, 21).Value = 2
$ws.Cells.Item(28, 22).Value = 2

# Row 29
$ws.Cells.Item(29, 9).Value = 'target'
$ws.Cells.Item(29, 10).Value = 'old'
$ws.Cells.Item(29, 11).Value = 'j'
$ws.Cells.Item(29, 12).Value = 'stimuli/img_ewrjk.png'
$ws.Cells.Item(29, 13).Value = 73.09090909090909
$ws.Cells.Item(29, 14).Value = 53.39393939393939
$ws.Cells.Item(29, 15).Value = 63.24242424242424
$ws.Cells.Item(29, 16).Value = 33
$ws.Cells.Item(29, 17).Value = 6
$ws.Cells.Item(29, 18).Value = 6
$ws.Cells.Item(29, 19).Value = 6
$ws.Cells.Item(29, 20).Value = 6
$ws.Cells.Item(29, 21).Value = 6
$ws.Cells.Item(29, 22).Value = 6

# Row 30
$ws.Cells.Item(30, 12).Value = 'stimuli/img_lszzj.png'
$ws.Cells.Item(30, 13).Value = 64.70588235294117
$ws.Cells.Item(30, 14).Value = 45.58823529411764
$ws.Cells.Item(30, 15).Value = 55.14705882352941
$ws.Cells.Item(30, 16).Value = 34
$ws.Cells.Item(30, 17).Value = 4
$ws.Cells.Item(30, 18).Value = 4
$ws.Cells.Item(30, 19).Value = 4
$ws.Cells.Item(30, 20).Value = 4
$ws.Cells.Item(30, 21).Value = 4
$ws.Cells.Item(30, 22).Value = 4

# Row 31
$ws.Cells.Item(31, 8).ClearContents()
$ws.Cells.Item(31, 9).ClearContents()
$ws.Cells.Item(31, 10).Value = 'catch'
$ws.Cells.Item(31, 11).Value = 'f'
$ws.Cells.Item(31, 12).Value = 'stimuli/catch_20.jpg'
$ws.Cells.Item(31, 13).ClearContents()
$ws.Cells.Item(31, 14).ClearContents()
$ws.Cells.Item(31, 15).ClearContents()
$ws.Cells.Item(31, 16).ClearContents()
$ws.Cells.Item(31, 17).ClearContents()
$ws.Cells.Item(31, 18).ClearContents()
$ws.Cells.Item(31, 19).ClearContents()
$ws.Cells.Item(31, 20).ClearContents()
$ws.Cells.Item(31, 21).ClearContents()
$ws.Cells.Item(31, 22).ClearContents()

# Row 32
$ws.Cells.Item(32, 9).ClearContents()
$ws.Cells.Item(32, 10).Value = 'new'
$ws.Cells.Item(32, 11).Value = 'f'
$ws.Cells.Item(32, 12).Value = 'stimuli/img_cv6mf.png'
$ws.Cells.Item(32, 13).Value = 66.8
$ws.Cells.Item(32, 14).Value = 42.08
$ws.Cells.Item(32, 15).Value = 54.44
$ws.Cells.Item(32, 16).Value = 25
$ws.Cells.Item(32, 17).Value = 4
$ws.Cells.Item(32, 18).Value = 4
$ws.Cells.Item(32, 19).Value = 4
$ws.Cells.Item(32, 20).Value = 4
$ws.Cells.Item(32, 21).Value = 4
$ws.Cells.Item(32, 22).Value = 4

# Row 33
$ws.Cells.Item(33, 9).ClearContents()
$ws.Cells.Item(33, 10).Value = 'new'
$ws.Cells.Item(33, 11).Value = 'f'
$ws.Cells.Item(33, 12).Value = 'stimuli/img_yeh72.png'
$ws.Cells.Item(33, 13).Value = 68.66666666666667
$ws.Cells.Item(33, 14).Value = 45.21212121212121
$ws.Cells.Item(33, 15).Value = 56.93939393939394
$ws.Cells.Item(33, 16).Value = 33
$ws.Cells.Item(33, 17).Value = 4
$ws.Cells.Item(33, 18).Value = 4
$ws.Cells.Item(33, 19).Value = 4
$ws.Cells.Item(33, 20).Value = 4
$ws.Cells.Item(33, 21).Value = 4
$ws.Cells.Item(33, 22).Value = 4

# Row 34
$ws.Cells.Item(34, 9).Value = 'target'
$ws.Cells.Item(34, 10).Value = 'old'
$ws.Cells.Item(34, 11).Value = 'j'
$ws.Cells.Item(34, 12).Value = 'stimuli/img_mjxmq.png'
$ws.Cells.Item(34, 13).Value = 77.07692307692308
$ws.Cells.Item(34, 14).Value = 58.15384615384615
$ws.Cells.Item(34, 15).Value = 67.61538461538461
$ws.Cells.Item(34, 16).Value = 39
$ws.Cells.Item(34, 17).Value = 7
$ws.Cells.Item(34, 18).Value = 7
$ws.Cells.Item(34, 19).Value = 7
$ws.Cells.Item(34, 20).Value = 7
$ws.Cells.Item(34, 21).Value = 7
$ws.Cells.Item(34, 22).Value = 7

# Row 35
$ws.Cells.Item(35, 9).Value = 'target'
$ws.Cells.Item(35, 10).Value = 'old'
$ws.Cells.Item(35, 11).Value = 'j'
$ws.Cells.Item(35, 12).Value = 'stimuli/img_vbrb7.png'
$ws.Cells.Item(35, 13).Value = 85.5625
$ws.Cells.Item(35, 14).Value = 71.46875
$ws.Cells.Item(35, 15).Value = 78.515625
$ws.Cells.Item(35, 16).Value = 32
$ws.Cells.Item(35, 17).Value = 10
$ws.Cells.Item(35, 18).Value = 10
$ws.Cells.Item(35, 19).Value = 10
$ws.Cells.Item(35, 20).Value = 10
$ws.Cells.Item(35, 21).Value = 10
$ws.Cells.Item(35, 22).Value = 10

# Row 36
$ws.Cells.Item(36, 12).Value = 'stimuli/img_t90e2.png'
$ws.Cells.Item(36, 13).Value = 83.0625
$ws.Cells.Item(36, 14).Value = 61.96875
$ws.Cells.Item(36, 15).Value = 72.515625
$ws.Cells.Item(36, 16).Value = 32
$ws.Cells.Item(36, 17).Value = 9
$ws.Cells.Item(36, 18).Value = 9
$ws.Cells.Item(36, 19).Value = 9
$ws.Cells.Item(36, 20).Value = 9
$ws.Cells.Item(36, 21).Value = 9
$ws.Cells.Item(36, 22).Value = 9

# Row 37
$ws.Cells.Item(37, 12).Value = 'stimuli/img_p659z.png'
$ws.Cells.Item(37, 13).Value = 84.21621621621621
$ws.Cells.Item(37, 14).Value = 65.37837837837837
$ws.Cells.Item(37, 15).Value = 74.79729729729729
$ws.Cells.Item(37, 16).Value = 37

# Row 38
$ws.Cells.Item(38, 12).Value = 'stimuli/img_r2lxk.png'
$ws.Cells.Item(38, 13).Value = 89.24242424242425
$ws.Cells.Item(38, 14).Value = 67.6969696969697
$ws.Cells.Item(38, 15).Value = 78.46969696969697
$ws.Cells.Item(38, 17).Value = 10
$ws.Cells.Item(38, 18).Value = 10
$ws.Cells.Item(38, 19).Value = 10
$ws.Cells.Item(38, 20).Value = 10
$ws.Cells.Item(38, 21).Value = 10
$ws.Cells.Item(38, 22).Value = 10

# Row 39
$ws.Cells.Item(39, 12).Value = 'stimuli/img_7ed9m.png'
$ws.Cells.Item(39, 13).Value = 80.71875
$ws.Cells.Item(39, 14).Value = 58.65625
$ws.Cells.Item(39, 15).Value = 69.6875
$ws.Cells.Item(39, 16).Value = 32
$ws.Cells.Item(39, 17).Value = 8
$ws.Cells.Item(39, 18).Value = 8
$ws.Cells.Item(39, 19).Value = 8
$ws.Cells.Item(39, 20).Value = 8
$ws.Cells.Item(39, 21).Value = 8
$ws.Cells.Item(39, 22).Value = 8

# Row 40
$ws.Cells.Item(40, 9).Value = 'target'
$ws.Cells.Item(40, 10).Value = 'old'
$ws.Cells.Item(40, 11).Value = 'j'
$ws.Cells.Item(40, 12).Value = 'stimuli/img_z293c.png'
$ws.Cells.Item(40, 13).Value = 71.26470588235294
$ws.Cells.Item(40, 14).Value = 46.88235294117647
$ws.Cells.Item(40, 15).Value = 59.07352941176471
$ws.Cells.Item(40, 16).Value = 34
$ws.Cells.Item(40, 17).Value = 5
$ws.Cells.Item(40, 18).Value = 5
$ws.Cells.Item(40, 19).Value = 5
$ws.Cells.Item(40, 20).Value = 5
$ws.Cells.Item(40, 21).Value = 5
$ws.Cells.Item(40, 22).Value = 5

# Row 41
$ws.Cells.Item(41, 12).Value = 'stimuli/img_7wul8.png'
$ws.Cells.Item(41, 13).Value = 43.03030303030303
$ws.Cells.Item(41, 14).Value = 25.54545454545455
$ws.Cells.Item(41, 15).Value = 34.28787878787879
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = 1
$ws.Cells.Item(41, 19).Value = 1
$ws.Cells.Item(41, 20).Value = 1
$ws.Cells.Item(41, 21).Value = 1
$ws.Cells.Item(41, 22).Value = 1

# Row 42
$ws.Cells.Item(42, 12).Value = 'stimuli/img_njmgp.png'
$ws.Cells.Item(42, 13).Value = 80.48148148148148
$ws.Cells.Item(42, 14).Value = 58.4074074074074
$ws.Cells.Item(42, 15).Value = 69.44444444444444
$ws.Cells.Item(42, 16).Value = 27
$ws.Cells.Item(42, 17).Value = 8
$ws.Cells.Item(42, 18).Value = 8
$ws.Cells.Item(42, 19).Value = 8
$ws.Cells.Item(42, 20).Value = 8
$ws.Cells.Item(42, 21).Value = 8
$ws.Cells.Item(42, 22).Value = 8
